# Inserts one new weekly price record for "Brócoli" (Terminal Hortofrutícola
# Agro Chillán) into the data table on Sheet1, at row 364 — pushing every
# existing row from 364 down through the former last row (382) down by one
# (new rows 365..383), growing the used range from A1:R382 to A1:R383.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 364 (shifts 364:382 -> 365:383).
# Excel's native row-insert also carries the formatting of the row above
# into the new row, so column D picks up the existing date number format.
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new record's data.
$ws.Cells.Item(364, 1).Value  = 7
$ws.Cells.Item(364, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(364, 3).Value  = "Ñuble"
$ws.Cells.Item(364, 4).Value  = 44939
$ws.Cells.Item(364, 5).Value  = 16
$ws.Cells.Item(364, 6).Value  = 100112023
$ws.Cells.Item(364, 7).Value  = "Brócoli"
$ws.Cells.Item(364, 8).Value  = "Sin especificar"
$ws.Cells.Item(364, 9).Value  = "Primera"
$ws.Cells.Item(364, 10).Value = 200
$ws.Cells.Item(364, 11).Value = 700
$ws.Cells.Item(364, 12).Value = 750
$ws.Cells.Item(364, 13).Value = 725
$ws.Cells.Item(364, 14).Value = "$/unidad"
$ws.Cells.Item(364, 15).Value = "Región del Maule"
$ws.Cells.Item(364, 16).Value = 725
$ws.Cells.Item(364, 17).Value = 1
$ws.Cells.Item(364, 18).Value = "Hortaliza"
